$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Simple text replacements (Find & Replace on exact phrases)
# ------------------------------------------------------------------

$d.Content.Find.Execute(
    "Senior Software Engineer with 21 years of experience in full-stack development, data engineering, and scalable web applications. Expert in Python, Drupal (4-10), GeoDjango, Flask, and cloud architecture with proven track record building enterprise-scale systems.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Senior Software Engineer with 15+ years of experience in full-stack development, data engineering, and scalable web applications. Expert in Python, Drupal (4-10), GeoDjango, Flask, and cloud architecture with proven track record building enterprise-scale systems.",
    2) | Out-Null

$d.Content.Find.Execute(
    "• Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys",
    2) | Out-Null

$d.Content.Find.Execute(
    "• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously",
    2) | Out-Null

$d.Content.Find.Execute(
    "• Integrated mapping and visualization tools for political campaign data analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "• Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs",
    2) | Out-Null

# ------------------------------------------------------------------
# 2) Structural edits (insert / delete whole paragraphs).
#    Done from the bottom of the document upward so that earlier
#    paragraph indices are not shifted by later insertions/deletions.
# ------------------------------------------------------------------

# --- EDUCATION section: remove the heading + two degree lines, and
#     replace them with a new bullet in the Feldman Group block ---
$pEduStart = $d.Paragraphs.Item(73)
$pEduEnd = $d.Paragraphs.Item(75)
$rngEdu = $d.Range($pEduStart.Range.Start, $pEduEnd.Range.End)
$rngEdu.Delete()

$pFeldmanLastBullet = $d.Paragraphs.Item(72)
$pFeldmanLastBullet.Range.InsertParagraphAfter()
$d.Paragraphs.Item(73).Range.Text = "• Trained staff on PHP/MySQL for data analysis and reporting systems"

# --- Lake Research Partners: add Python tooling training bullet ---
$pLakeLastBullet = $d.Paragraphs.Item(64)
$pLakeLastBullet.Range.InsertParagraphAfter()
$d.Paragraphs.Item(65).Range.Text = "• Trained staff on building Python tooling for report generation and analysis"

# --- Praxis Project: add Drupal sites bullet ---
$pPraxisLastBullet = $d.Paragraphs.Item(56)
$pPraxisLastBullet.Range.InsertParagraphAfter()
$d.Paragraphs.Item(57).Range.Text = "• Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation"

# --- Salsa Labs: add billions-of-records bullet ---
$pSalsaLastBullet = $d.Paragraphs.Item(48)
$pSalsaLastBullet.Range.InsertParagraphAfter()
$d.Paragraphs.Item(49).Range.Text = "• Handled billions of records with millions of columns in high-performance CRM system"

Write-Output "Edit complete"
